$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 6
$ws.Range("D6").Value = 6
$ws.Range("E6").Value = 4
$ws.Range("A15").Value = "Ajans Of"
$ws.Range("B15").Value = "Fortuna United"
$ws.Range("E19").Select()
